$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.733.04"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").Value = "1.862.11"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("E4").Value = "  -0.92%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.020"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4379"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3808"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07447"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8833"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.63"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("D12").Value = "1.865.84"
$ws.Range("E12").Value = "  +0.86%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.754"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.501"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07143"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.025"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009091"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.020"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.50"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").Value = "27.734.33"
$ws.Range("E21").Value = "  +0.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.292"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.53%  "

$ws.Range("D24").Value = "2.089.91"
$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.048"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.50"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.005"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.373"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.50"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09057"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.213"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7669"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.034"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.568"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.021"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.142"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01981"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05297"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.873"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5195"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.950"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1681"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.710"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.78"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "110.06"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.717"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.021"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.83%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06503"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.61%  "

$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4718"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.884"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.85%  "
